$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 266, shifting existing rows 266-337 down to 267-338.
$ws.Rows("266:266").Insert()

# Populate the newly inserted row 266 with the new weekly data point.
$ws.Range("A266").Value2 = 3
$ws.Range("B266").Value2 = "Femacal de La Calera"
$ws.Range("C266").Value2 = "Coquimbo"
$ws.Range("D266").Value2 = 44642
$ws.Range("D266").NumberFormat = $ws.Range("D267").NumberFormat
$ws.Range("E266").Value2 = 5
$ws.Range("F266").Value2 = 100112031
$ws.Range("G266").Value2 = "Poroto verde"
$ws.Range("H266").Value2 = "Magnum"
$ws.Range("I266").Value2 = "Primera"
$ws.Range("J266").Value2 = 73
$ws.Range("K266").Value2 = 21000
$ws.Range("L266").Value2 = 22000
$ws.Range("M266").Value2 = 21479
$ws.Range("N266").Value2 = "$/malla 25 kilos"
$ws.Range("O266").Value2 = "Provincia de Santiago"
$ws.Range("P266").Value2 = 859
$ws.Range("Q266").Value2 = 25
$ws.Range("R266").Value2 = "Hortaliza"
